$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 578, shifting existing rows 578..619 down to 580..621
$ws.Rows.Item(578).Resize(2).Insert()

# Fill in the two newly inserted rows with the new data.
# Force text formatting first so "2026/01/08" is stored as a literal
# string (matching the sheet's existing inline-string date cells)
# instead of being auto-parsed into a date serial, then reset the
# style back to Normal so no stray style index is left behind.
$ws.Range("A578:A579").NumberFormat = "@"

$ws.Range("A578").Value = "2026/01/08"
$ws.Range("B578").Value = "木"
$ws.Range("C578").Value = 7
$ws.Range("D578").Value = 25

$ws.Range("A579").Value = "2026/01/08"
$ws.Range("B579").Value = "木"
$ws.Range("C579").Value = 10
$ws.Range("D579").Value = 25

$ws.Range("A578:D579").Style = "Normal"
